# Clear unused code in IAM Module.
# Updates the "Test Cases" sheet (IAM016 row) Jira id / description text,
# appends a new IAM035 test-case row, and updates view/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Row 17 (IAM016 - change password test case): extend Jira id list and
#     description with the new password-reset-workflow requirement.
$ws.Cells.Item(17, 3).Value = "Verify that user is able to change his STeAM password by using Forgot Password link and that he is able to login with his new password||Verify Upon completion of establishing a new password, a user who wants to go to Neon shall be presented a confirmation page with an optional link back to Neon Landing page.||Verify that,the system should support a Neon password reset workflow with the following configurations:"
$ws.Cells.Item(17, 2).Value = "OPQA-535||OPQA-1955||OPQA-3686"
$ws.Rows.Item(17).RowHeight = 86.4

# --- New row 36: IAM035 test case (email leading/trailing space trimming).
$ws.Range("A35:E35").Copy()
$ws.Range("A36:E36").PasteSpecial(-4122)
$ws.Rows.Item(36).RowHeight = 28.8
$ws.Cells.Item(36, 1).Value = "IAM035"
$ws.Cells.Item(36, 2).Value = "OPQA-1851"
$ws.Cells.Item(36, 3).Value = "Verify that system should remove any leading or trailing spaces of an email address entered by the user before validating it."
$ws.Cells.Item(36, 4).Value = "Y"
$ws.Cells.Item(36, 5).Value = "PASS"

# --- View state: IAM019 sheet selection moves to C16, then back to the
#     "Test Cases" sheet (kept active) with the new row selected.
$ws19 = $wb.Worksheets.Item("IAM019")
$ws19.Range("C16").Select()

$ws.Activate()
$ws.Range("C36").Select()
